$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '64.706.38'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '3.371.47'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.79'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.34%  '
$ws.Range('E7').Value = '  +2.08%  '
$ws.Range('D8').Value = '3.362.25'
$ws.Range('E8').Value = '  +1.73%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('E10').Value = '  +8.35%  '
$ws.Range('E11').Value = '  +3.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.85'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('E13').Value = '  +3.93%  '
$ws.Range('E14').Value = '  +2.28%  '
$ws.Range('D15').Value = '3.904.59'
$ws.Range('E15').Value = '  +1.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.23'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('E17').Value = '  +2.17%  '
$ws.Range('D18').Value = '3.377.29'
$ws.Range('E18').Value = '  +1.82%  '
$ws.Range('D19').Value = '64.673.35'
$ws.Range('E19').Value = '  +2.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.990'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '459.26'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.25%  '
$ws.Range('E23').Value = '  +8.78%  '
$ws.Range('E24').Value = '  +1.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.27'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.58'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.65%  '
$ws.Range('E28').Value = '  +1.78%  '
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.60'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.67'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.84%  '
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '571.31'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '61.12'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.12%  '
$ws.Range('E35').Value = '  +1.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.64'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.78%  '
$ws.Range('E38').Value = '  -4.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.39'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('D40').Value = '0.0₃0741'
$ws.Range('E40').Value = '  +0.37%  '
$ws.Range('E41').Value = '  +1.52%  '
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('D43').Value = '3.074.82'
$ws.Range('E43').Value = '  -1.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.82'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0415'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.52%  '
$ws.Range('E46').Value = '  +4.66%  '
$ws.Range('E47').Value = '  +1.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.11'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.34%  '
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '138.27'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.20'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.20%  '
